$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "normal update`n( 1. Find the Robot using the Unit Name, machine name and user name`n  2. check if the robot name is not the same as to-be robot name already`n  3. update the robot name to the new name`n  4. find the robot using the Unit name, machine name and user name again`n  5. confirm the robot name is the same as to-be robot name.)"

$ws.Range("A9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 129.6

$ws.Range("A9").Select()
